$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 11 ("Sub total" line) moves down to row 15, making room for
# two new component rows (11 and 13) each followed by a blank spacer row
# (12 and 14), matching the existing blank-spacer-then-total layout.

# Clear out the old "Sub total" row content (it will be re-created at row 15).
$ws.Range("O11").ClearContents() | Out-Null
$ws.Range("P11").ClearContents() | Out-Null

# Give the new spacer/blank column-D cells in rows 11-15 the same format
# used by the existing blank spacer at D10.
$ws.Range("D10").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4104) | Out-Null
$ws.Range("D12").PasteSpecial(-4104) | Out-Null
$ws.Range("D13").PasteSpecial(-4104) | Out-Null
$ws.Range("D14").PasteSpecial(-4104) | Out-Null
$ws.Range("D15").PasteSpecial(-4104) | Out-Null

# New item row 11: PCB
$ws.Range("C11").Value = "PCB"
$ws.Range("E11").Value = 1
$ws.Range("P9").Copy() | Out-Null
$ws.Range("P11").PasteSpecial(-4104) | Out-Null
$ws.Range("P11").Value = 0.5

# New item row 13: Assembly
$ws.Range("C13").Value = "Assembly"
$ws.Range("E13").Value = 1
$ws.Range("P9").Copy() | Out-Null
$ws.Range("P13").PasteSpecial(-4104) | Out-Null
$ws.Range("P13").Value = 0.5

# "Sub total" label + formula now live on row 15 and sum through the new rows.
$ws.Range("O15").Value = "Sub total"
$ws.Range("P9").Copy() | Out-Null
$ws.Range("P15").PasteSpecial(-4104) | Out-Null
$ws.Range("P15").Formula = "=SUM(P2:P13)"

$excel.CutCopyMode = 0
$ws.Range("P15").Select()
